$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-09-25 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-09-26 Thursday", 2) | Out-Null
$d.Content.Find.Execute("88-40=48", $true, $false, $false, $false, $false, $true, 1, $false, "87-73=14", 2) | Out-Null
$d.Content.Find.Execute("39-9=30", $true, $false, $false, $false, $false, $true, 1, $false, "39-33=6", 2) | Out-Null
$d.Content.Find.Execute("55+44=99", $true, $false, $false, $false, $false, $true, 1, $false, "32+25=57", 2) | Out-Null
$d.Content.Find.Execute("81-13=68", $true, $false, $false, $false, $false, $true, 1, $false, "76-7=69", 2) | Out-Null
$d.Content.Find.Execute("71+3=74", $true, $false, $false, $false, $false, $true, 1, $false, "32-24=8", 2) | Out-Null
$d.Content.Find.Execute("61+7=68", $true, $false, $false, $false, $false, $true, 1, $false, "78-33=45", 2) | Out-Null
$d.Content.Find.Execute("3+82=85", $true, $false, $false, $false, $false, $true, 1, $false, "45+43=88", 2) | Out-Null
$d.Content.Find.Execute("27+23=50", $true, $false, $false, $false, $false, $true, 1, $false, "50+9=59", 2) | Out-Null
$d.Content.Find.Execute("73+8=81", $true, $false, $false, $false, $false, $true, 1, $false, "33+65=98", 2) | Out-Null
$d.Content.Find.Execute("78+20=98", $true, $false, $false, $false, $false, $true, 1, $false, "70-59=11", 2) | Out-Null
$d.Content.Find.Execute("57-15=42", $true, $false, $false, $false, $false, $true, 1, $false, "49+15=64", 2) | Out-Null
$d.Content.Find.Execute("13-7=6", $true, $false, $false, $false, $false, $true, 1, $false, "11+17=28", 2) | Out-Null
$d.Content.Find.Execute("30+23=53", $true, $false, $false, $false, $false, $true, 1, $false, "36-16=20", 2) | Out-Null
$d.Content.Find.Execute("31-0=31", $true, $false, $false, $false, $false, $true, 1, $false, "0+40=40", 2) | Out-Null
$d.Content.Find.Execute("23+64=87", $true, $false, $false, $false, $false, $true, 1, $false, "14+13=27", 2) | Out-Null
$d.Content.Find.Execute("58+10=68", $true, $false, $false, $false, $false, $true, 1, $false, "34-18=16", 2) | Out-Null
$d.Content.Find.Execute("1+20=21", $true, $false, $false, $false, $false, $true, 1, $false, "29+26=55", 2) | Out-Null
$d.Content.Find.Execute("29+52=81", $true, $false, $false, $false, $false, $true, 1, $false, "5+66=71", 2) | Out-Null
$d.Content.Find.Execute("53+37=90", $true, $false, $false, $false, $false, $true, 1, $false, "71-68=3", 2) | Out-Null
$d.Content.Find.Execute("35+18=53", $true, $false, $false, $false, $false, $true, 1, $false, "19+7=26", 2) | Out-Null
$d.Content.Find.Execute("32+59=91", $true, $false, $false, $false, $false, $true, 1, $false, "48-44=4", 2) | Out-Null
$d.Content.Find.Execute("35+63=98", $true, $false, $false, $false, $false, $true, 1, $false, "36-5=31", 2) | Out-Null
$d.Content.Find.Execute("41+56=97", $true, $false, $false, $false, $false, $true, 1, $false, "38+9=47", 2) | Out-Null
$d.Content.Find.Execute("79-22=57", $true, $false, $false, $false, $false, $true, 1, $false, "56+15=71", 2) | Out-Null
$d.Content.Find.Execute("64+33=97", $true, $false, $false, $false, $false, $true, 1, $false, "23+12=35", 2) | Out-Null
$d.Content.Find.Execute("22+65=87", $true, $false, $false, $false, $false, $true, 1, $false, "69-5=64", 2) | Out-Null
$d.Content.Find.Execute("18-1=17", $true, $false, $false, $false, $false, $true, 1, $false, "35-8=27", 2) | Out-Null
$d.Content.Find.Execute("6+29=35", $true, $false, $false, $false, $false, $true, 1, $false, "39-17=22", 2) | Out-Null
$d.Content.Find.Execute("14-11=3", $true, $false, $false, $false, $false, $true, 1, $false, "75-43=32", 2) | Out-Null
$d.Content.Find.Execute("76-1=75", $true, $false, $false, $false, $false, $true, 1, $false, "86-47=39", 2) | Out-Null
$d.Content.Find.Execute("14+52=66", $true, $false, $false, $false, $false, $true, 1, $false, "88-51=37", 2) | Out-Null
$d.Content.Find.Execute("43-12=31", $true, $false, $false, $false, $false, $true, 1, $false, "38+5=43", 2) | Out-Null
$d.Content.Find.Execute("51+36=87", $true, $false, $false, $false, $false, $true, 1, $false, "32+17=49", 2) | Out-Null
$d.Content.Find.Execute("4+85=89", $true, $false, $false, $false, $false, $true, 1, $false, "93-79=14", 2) | Out-Null
$d.Content.Find.Execute("95-59=36", $true, $false, $false, $false, $false, $true, 1, $false, "83-78=5", 2) | Out-Null
$d.Content.Find.Execute("73-21=52", $true, $false, $false, $false, $false, $true, 1, $false, "95-37=58", 2) | Out-Null
$d.Content.Find.Execute("65+11=76", $true, $false, $false, $false, $false, $true, 1, $false, "67+5=72", 2) | Out-Null
$d.Content.Find.Execute("14+15=29", $true, $false, $false, $false, $false, $true, 1, $false, "38-21=17", 2) | Out-Null
$d.Content.Find.Execute("96-2=94", $true, $false, $false, $false, $false, $true, 1, $false, "22+46=68", 2) | Out-Null
$d.Content.Find.Execute("92-45=47", $true, $false, $false, $false, $false, $true, 1, $false, "61-41=20", 2) | Out-Null
$d.Content.Find.Execute("44+29=73", $true, $false, $false, $false, $false, $true, 1, $false, "82-15=67", 2) | Out-Null
$d.Content.Find.Execute("56+35=91", $true, $false, $false, $false, $false, $true, 1, $false, "14+44=58", 2) | Out-Null
$d.Content.Find.Execute("76+2=78", $true, $false, $false, $false, $false, $true, 1, $false, "69-64=5", 2) | Out-Null
$d.Content.Find.Execute("39-24=15", $true, $false, $false, $false, $false, $true, 1, $false, "44+5=49", 2) | Out-Null
$d.Content.Find.Execute("71-31=40", $true, $false, $false, $false, $false, $true, 1, $false, "46+53=99", 2) | Out-Null
$d.Content.Find.Execute("21+15=36", $true, $false, $false, $false, $false, $true, 1, $false, "45-21=24", 2) | Out-Null
$d.Content.Find.Execute("15+26=41", $true, $false, $false, $false, $false, $true, 1, $false, "55+29=84", 2) | Out-Null
$d.Content.Find.Execute("64-43=21", $true, $false, $false, $false, $false, $true, 1, $false, "63+15=78", 2) | Out-Null
$d.Content.Find.Execute("31+8=39", $true, $false, $false, $false, $false, $true, 1, $false, "25+9=34", 2) | Out-Null
$d.Content.Find.Execute("99-87=12", $true, $false, $false, $false, $false, $true, 1, $false, "86-35=51", 2) | Out-Null
$d.Content.Find.Execute("93+6=99", $true, $false, $false, $false, $false, $true, 1, $false, "66-14=52", 2) | Out-Null
$d.Content.Find.Execute("67+31=98", $true, $false, $false, $false, $false, $true, 1, $false, "79+9=88", 2) | Out-Null
$d.Content.Find.Execute("12+36=48", $true, $false, $false, $false, $false, $true, 1, $false, "89-79=10", 2) | Out-Null
$d.Content.Find.Execute("83-62=21", $true, $false, $false, $false, $false, $true, 1, $false, "31+33=64", 2) | Out-Null
$d.Content.Find.Execute("80-17=63", $true, $false, $false, $false, $false, $true, 1, $false, "32+19=51", 2) | Out-Null
$d.Content.Find.Execute("70-44=26", $true, $false, $false, $false, $false, $true, 1, $false, "79-56=23", 2) | Out-Null
$d.Content.Find.Execute("59-22=37", $true, $false, $false, $false, $false, $true, 1, $false, "79-1=78", 2) | Out-Null
$d.Content.Find.Execute("49-13=36", $true, $false, $false, $false, $false, $true, 1, $false, "77+14=91", 2) | Out-Null
$d.Content.Find.Execute("12+65=77", $true, $false, $false, $false, $false, $true, 1, $false, "83-9=74", 2) | Out-Null
$d.Content.Find.Execute("47+43=90", $true, $false, $false, $false, $false, $true, 1, $false, "42+47=89", 2) | Out-Null
$d.Content.Find.Execute("79-33=46", $true, $false, $false, $false, $false, $true, 1, $false, "37+5=42", 2) | Out-Null
$d.Content.Find.Execute("33-2=31", $true, $false, $false, $false, $false, $true, 1, $false, "6+66=72", 2) | Out-Null
$d.Content.Find.Execute("56-36=20", $true, $false, $false, $false, $false, $true, 1, $false, "17-1=16", 2) | Out-Null
$d.Content.Find.Execute("67-53=14", $true, $false, $false, $false, $false, $true, 1, $false, "3+84=87", 2) | Out-Null
$d.Content.Find.Execute("36-23=13", $true, $false, $false, $false, $false, $true, 1, $false, "35+16=51", 2) | Out-Null
$d.Content.Find.Execute("64+31=95", $true, $false, $false, $false, $false, $true, 1, $false, "85-56=29", 2) | Out-Null
$d.Content.Find.Execute("1+57=58", $true, $false, $false, $false, $false, $true, 1, $false, "74-59=15", 2) | Out-Null
$d.Content.Find.Execute("58-22=36", $true, $false, $false, $false, $false, $true, 1, $false, "10+88=98", 2) | Out-Null
$d.Content.Find.Execute("39-4=35", $true, $false, $false, $false, $false, $true, 1, $false, "36-31=5", 2) | Out-Null
$d.Content.Find.Execute("29-18=11", $true, $false, $false, $false, $false, $true, 1, $false, "88-68=20", 2) | Out-Null
$d.Content.Find.Execute("50-32=18", $true, $false, $false, $false, $false, $true, 1, $false, "34-27=7", 2) | Out-Null
$d.Content.Find.Execute("38-38=0", $true, $false, $false, $false, $false, $true, 1, $false, "80-34=46", 2) | Out-Null
$d.Content.Find.Execute("12+29=41", $true, $false, $false, $false, $false, $true, 1, $false, "45-19=26", 2) | Out-Null
$d.Content.Find.Execute("83-44=39", $true, $false, $false, $false, $false, $true, 1, $false, "51+43=94", 2) | Out-Null
$d.Content.Find.Execute("65-29=36", $true, $false, $false, $false, $false, $true, 1, $false, "47-47=0", 2) | Out-Null
$d.Content.Find.Execute("96-56=40", $true, $false, $false, $false, $false, $true, 1, $false, "98-86=12", 2) | Out-Null
$d.Content.Find.Execute("65-36=29", $true, $false, $false, $false, $false, $true, 1, $false, "81-50=31", 2) | Out-Null
$d.Content.Find.Execute("33-21=12", $true, $false, $false, $false, $false, $true, 1, $false, "83-21=62", 2) | Out-Null
$d.Content.Find.Execute("25+49=74", $true, $false, $false, $false, $false, $true, 1, $false, "87-25=62", 2) | Out-Null
$d.Content.Find.Execute("83-52=31", $true, $false, $false, $false, $false, $true, 1, $false, "82-68=14", 2) | Out-Null
$d.Content.Find.Execute("43+0=43", $true, $false, $false, $false, $false, $true, 1, $false, "26+1=27", 2) | Out-Null
$d.Content.Find.Execute("38-4=34", $true, $false, $false, $false, $false, $true, 1, $false, "40-27=13", 2) | Out-Null
$d.Content.Find.Execute("99-31=68", $true, $false, $false, $false, $false, $true, 1, $false, "92-28=64", 2) | Out-Null
$d.Content.Find.Execute("24-22=2", $true, $false, $false, $false, $false, $true, 1, $false, "38-35=3", 2) | Out-Null
$d.Content.Find.Execute("79-15=64", $true, $false, $false, $false, $false, $true, 1, $false, "30-5=25", 2) | Out-Null
$d.Content.Find.Execute("42+44=86", $true, $false, $false, $false, $false, $true, 1, $false, "2+94=96", 2) | Out-Null
$d.Content.Find.Execute("69+22=91", $true, $false, $false, $false, $false, $true, 1, $false, "32+12=44", 2) | Out-Null
$d.Content.Find.Execute("12+77=89", $true, $false, $false, $false, $false, $true, 1, $false, "66-19=47", 2) | Out-Null
$d.Content.Find.Execute("75-54=21", $true, $false, $false, $false, $false, $true, 1, $false, "85-18=67", 2) | Out-Null
$d.Content.Find.Execute("52-35=17", $true, $false, $false, $false, $false, $true, 1, $false, "13+51=64", 2) | Out-Null
$d.Content.Find.Execute("99-98=1", $true, $false, $false, $false, $false, $true, 1, $false, "65+32=97", 2) | Out-Null
$d.Content.Find.Execute("84-2=82", $true, $false, $false, $false, $false, $true, 1, $false, "41+52=93", 2) | Out-Null
$d.Content.Find.Execute("26+11=37", $true, $false, $false, $false, $false, $true, 1, $false, "63+12=75", 2) | Out-Null
$d.Content.Find.Execute("73+0=73", $true, $false, $false, $false, $false, $true, 1, $false, "56-44=12", 2) | Out-Null
$d.Content.Find.Execute("53-1=52", $true, $false, $false, $false, $false, $true, 1, $false, "46+20=66", 2) | Out-Null
$d.Content.Find.Execute("76-72=4", $true, $false, $false, $false, $false, $true, 1, $false, "47-19=28", 2) | Out-Null
$d.Content.Find.Execute("95-42=53", $true, $false, $false, $false, $false, $true, 1, $false, "57+24=81", 2) | Out-Null
$d.Content.Find.Execute("6+91=97", $true, $false, $false, $false, $false, $true, 1, $false, "81-40=41", 2) | Out-Null
$d.Content.Find.Execute("60+27=87", $true, $false, $false, $false, $false, $true, 1, $false, "21+25=46", 2) | Out-Null
$d.Content.Find.Execute("41-9=32", $true, $false, $false, $false, $false, $true, 1, $false, "9+72=81", 2) | Out-Null
